$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New region labels and sample counts for rows 2-25 (column C = region, column D = number_of_samples)
$ws.Range("C2").Value = "Cg25"
$ws.Range("D2").Value = 7

$ws.Range("C3").Value = "Cg25"
$ws.Range("D3").Value = 10

$ws.Range("C4").Value = "Nac"
$ws.Range("D4").Value = 9

$ws.Range("C5").Value = "Nac"
$ws.Range("D5").Value = 13

$ws.Range("C6").Value = "OFC"
$ws.Range("D6").Value = 9

$ws.Range("C7").Value = "OFC"
$ws.Range("D7").Value = 12

$ws.Range("C8").Value = "Sub"
$ws.Range("D8").Value = 7

$ws.Range("C9").Value = "Sub"
$ws.Range("D9").Value = 12

$ws.Range("C10").Value = "aINS"
$ws.Range("D10").Value = 9

$ws.Range("C11").Value = "aINS"
$ws.Range("D11").Value = 13

$ws.Range("C12").Value = "dlPFC"
$ws.Range("D12").Value = 9

$ws.Range("C13").Value = "dlPFC"
$ws.Range("D13").Value = 13

$ws.Range("C14").Value = "Cg25"
$ws.Range("D14").Value = 8

$ws.Range("C15").Value = "Cg25"
$ws.Range("D15").Value = 3

$ws.Range("C16").Value = "Nac"
$ws.Range("D16").Value = 13

$ws.Range("C17").Value = "Nac"
$ws.Range("D17").Value = 15

$ws.Range("C18").Value = "OFC"
$ws.Range("D18").Value = 13

$ws.Range("C19").Value = "OFC"
$ws.Range("D19").Value = 13

$ws.Range("C20").Value = "Sub"
$ws.Range("D20").Value = 12

$ws.Range("C21").Value = "Sub"
$ws.Range("D21").Value = 12

$ws.Range("C22").Value = "aINS"
$ws.Range("D22").Value = 13

$ws.Range("C23").Value = "aINS"
$ws.Range("D23").Value = 12

$ws.Range("C24").Value = "dlPFC"
$ws.Range("D24").Value = 13

$ws.Range("C25").Value = "dlPFC"
$ws.Range("D25").Value = 13
